$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Jalen Brunson'
$ws.Cells.Item(2, 2).Value = 'PG'
$ws.Cells.Item(2, 3).Value = 'New York Knicks'

$ws.Cells.Item(3, 1).Value = 'Coby White'
$ws.Cells.Item(3, 2).Value = 'PG,SG'
$ws.Cells.Item(3, 3).Value = 'Chicago Bulls'

$ws.Cells.Item(4, 1).Value = 'Trae Young'
$ws.Cells.Item(4, 2).Value = 'PG'
$ws.Cells.Item(4, 3).Value = 'Atlanta Hawks'

$ws.Cells.Item(5, 1).Value = 'Devin Booker'
$ws.Cells.Item(5, 2).Value = 'PG,SG'
$ws.Cells.Item(5, 3).Value = 'Phoenix Suns'

$ws.Cells.Item(6, 1).Value = 'Devin Vassell'
$ws.Cells.Item(6, 2).Value = 'SG,SF'
$ws.Cells.Item(6, 3).Value = 'San Antonio Spurs'

$ws.Cells.Item(7, 1).Value = 'LeBron James'
$ws.Cells.Item(7, 2).Value = 'SF,PF'
$ws.Cells.Item(7, 3).Value = 'Los Angeles Lakers'

$ws.Cells.Item(8, 1).Value = 'Norman Powell'
$ws.Cells.Item(8, 2).Value = 'SG,SF'
$ws.Cells.Item(8, 3).Value = 'LA Clippers'

$ws.Cells.Item(9, 1).Value = 'Santi Aldama'
$ws.Cells.Item(9, 2).Value = 'PF,C'
$ws.Cells.Item(9, 3).Value = 'Memphis Grizzlies'

$ws.Cells.Item(10, 1).Value = 'Alperen Sengün'
$ws.Cells.Item(10, 2).Value = 'C'
$ws.Cells.Item(10, 3).Value = 'Houston Rockets'

$ws.Cells.Item(11, 1).Value = 'Goga Bitadze'
$ws.Cells.Item(11, 2).Value = 'C'
$ws.Cells.Item(11, 3).Value = 'Orlando Magic'

$ws.Cells.Item(12, 1).Value = 'Desmond Bane'
$ws.Cells.Item(12, 2).Value = 'SG,SF'
$ws.Cells.Item(12, 3).Value = 'Memphis Grizzlies'

$ws.Cells.Item(13, 1).Value = 'Jared McCain'
$ws.Cells.Item(13, 2).Value = 'PG'
$ws.Cells.Item(13, 3).Value = 'Philadelphia 76ers'

$ws.Cells.Item(14, 1).Value = 'Jordan Clarkson'
$ws.Cells.Item(14, 2).Value = 'SG,SF'
$ws.Cells.Item(14, 3).Value = 'Utah Jazz'

$ws.Cells.Item(15, 1).Value = 'AJ Green'
$ws.Cells.Item(15, 2).Value = 'PG,SG'
$ws.Cells.Item(15, 3).Value = 'Milwaukee Bucks'

$ws.Cells.Item(16, 1).Value = 'Nicolas Claxton'
$ws.Cells.Item(16, 2).Value = 'C'
$ws.Cells.Item(16, 3).Value = 'Brooklyn Nets'

$ws.Cells.Item(17, 1).Value = 'Kawhi Leonard'
$ws.Cells.Item(17, 2).Value = 'SG,SF,PF'
$ws.Cells.Item(17, 3).Value = 'LA Clippers'

$ws.Cells.Item(18, 1).Value = 'Immanuel Quickley'
$ws.Cells.Item(18, 2).Value = 'PG,SG'
$ws.Cells.Item(18, 3).Value = 'Toronto Raptors'

$ws.Cells.Item(19, 1).Value = 'Walker Kessler'
$ws.Cells.Item(19, 2).Value = 'C'
$ws.Cells.Item(19, 3).Value = 'Utah Jazz'
